$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.914.84'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.799.13'
$ws.Range('E3').Value = '  +3.12%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '621.20'
$ws.Range('E5').Value = '  +3.74%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '177.35'
$ws.Range('E6').Value = '  -3.67%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.797.28'
$ws.Range('E7').Value = '  +3.13%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('E11').Value = '  -4.98%  '
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '40.66'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000262'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.439.67'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.807.24'
$ws.Range('E16').Value = '  +3.30%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '70.013.41'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.58'
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.86'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '509.91'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.57'
$ws.Range('E22').Value = '  +4.15%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.729'
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.52'
$ws.Range('E24').Value = '  +3.57%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '87.79'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '13.23'
$ws.Range('E27').Value = '  +29.08%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.07'
$ws.Range('E28').Value = '  +1.34%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -1.63%  '
$ws.Range('E32').Value = '  -4.58%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '31.48'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  +5.84%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.23'
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('E38').Value = '  +4.89%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.333'
$ws.Range('E39').Value = '  -2.90%  '
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '51.08'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '45.45'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '419.55'
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.039.02'
$ws.Range('E46').Value = '  -4.06%  '
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '27.36'
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '137.91'
$ws.Range('E50').Value = '  +0.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.47'
$ws.Range('E51').Value = '  +1.50%  '
